# Regenerate the experiment task-order sheets (new run of the order
# generation script). Each worksheet keeps its original rId/tab position,
# but is renamed and its B-column "task_order" values (a freshly generated
# set of stim-file orders) are rewritten. A couple of sheets grow/shrink by
# a row or two, so new rows are created via Range.Copy(destination) (which
# preserves the source cell's style, unlike a plain .Value assignment)
# before the new value is written on top, and extra rows are removed with
# EntireRow.Delete().

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# rId1 / sheet1.xml : GNG -> TOL  (5 data rows -> 7 data rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "TOL_TO-16515890392967472"
$ws.Range("A4").Copy($ws.Range("A6"))
$ws.Range("A4").Copy($ws.Range("A7"))

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "MM_stims-1651589039262676.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ZM_stims-16515890392470882.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "MM_stims-1651589039278299.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "ZM_stims-1651589039262676.csv"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "MM_stims-16515890392967472.csv"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "ZM_stims-1651589039278299.csv"

# ---------------------------------------------------------------------
# rId2 / sheet2.xml : NB -> NB  (9 data rows -> 9 data rows, values only)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16515890413674967"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "OB-16515890404340498.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "TB-16515890412217212.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "ZB-match_1-16515890403202424.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "ZB-match_1-16515890401731944.csv"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "OB-16515890408112302.csv"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "ZB-match_4-16515890393970292.csv"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "OB-16515890407955399.csv"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "TB-16515890412862198.csv"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "TB-16515890413496265.csv"

# ---------------------------------------------------------------------
# rId3 / sheet3.xml : RS -> vSAT  (2 data rows -> 4 data rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "vSAT_TO-16515890414299066"
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A3").Copy($ws.Range("A5"))

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "SAT_stims-16515890413674967.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "SAT_stims-16515890413830378.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "vSAT_stims-1651589041414283.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "vSAT_stims-16515890413986592.csv"

# ---------------------------------------------------------------------
# rId4 / sheet4.xml : TOL -> GNG  (7 data rows -> 5 data rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "GNG_TO-16515890414611573"
$ws.Rows.Item(7).EntireRow.Delete()
$ws.Rows.Item(6).EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "go_stims-16515890414299066.csv"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "GNG_stims-16515890414455323.csv"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "go_stims-16515890414455323.csv"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "GNG_stims-16515890414611573.csv"

# ---------------------------------------------------------------------
# rId5 / sheet5.xml : vSAT -> RS  (5 data rows -> 3 data rows)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Name = "RS_TO-16515890414611573"
$ws.Rows.Item(5).EntireRow.Delete()
$ws.Rows.Item(4).EntireRow.Delete()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eyes closed"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eyes open"
